$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 164, pushing the existing row 164 down to row 165.
$ws.Rows(164).Insert()

# New row 164 gets the data that row 163 used to hold (before today's update).
$ws.Range("A164").Value = 5
$ws.Range("B164").Value = "Macroferia Regional de Talca"
$ws.Range("C164").Value = "Maule"
$ws.Range("D164").Value = 44249
$ws.Range("D164").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E164").Value = 7
$ws.Range("F164").Value = 100112030
$ws.Range("G164").Value = "Poroto granado"
$ws.Range("H164").Value = "Sin especificar"
$ws.Range("I164").Value = "Primera"
$ws.Range("J164").Value = 300
$ws.Range("K164").Value = 25000
$ws.Range("L164").Value = 25000
$ws.Range("M164").Value = 25000
$ws.Range("N164").Value = "$/saco 25 kilos"
$ws.Range("O164").Value = "Región del Maule"
$ws.Range("P164").Value = 1000
$ws.Range("Q164").Value = 25
$ws.Range("R164").Value = "Hortaliza"

# Update row 163 with the new reported price data.
$ws.Range("D163").Value = 44939
$ws.Range("J163").Value = 500
$ws.Range("K163").Value = 38000
$ws.Range("L163").Value = 38000
$ws.Range("M163").Value = 38000
$ws.Range("P163").Value = 1520
